$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Galactus"
$ws.Range("B6").Value = "e7cf3ef4f17c3999a94f2c6f612e8a888e5b1026878e4e19398b23bd38ec221a"
$ws.Range("C6").Value = "galactus@gmail.com"
$ws.Range("D6").Value = "member"

$ws.Range("A7").Value = "Bowser"
$ws.Range("B7").Value = "e7cf3ef4f17c3999a94f2c6f612e8a888e5b1026878e4e19398b23bd38ec221a"
$ws.Range("C7").Value = "Bowser@Bowser.gov"
$ws.Range("D7").Value = "member"

$ws.Range("A8").Value = "BillGates"
$ws.Range("B8").Value = "e7cf3ef4f17c3999a94f2c6f612e8a888e5b1026878e4e19398b23bd38ec221a"
$ws.Range("C8").Value = "BillGates@Microsoft.com"
$ws.Range("D8").Value = "member"
